# Reassign the text labels of several tied rows in the goods-frequency
# table. Column A holds the good's name, column B its count. Re-running
# the generation script produced a different (but count-stable) ordering
# of labels among rows that share the same count value in column B.
# Only the labels in column A change here; the counts in column B are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    19 = "небогатый товар"
    20 = "крамными товар"
    21 = "мясо"
    22 = "железный товар"
    24 = "набойчатый товар"
    27 = "пушной товар"
    28 = "внутренний товар"
    30 = "питейный припасы"
    31 = "суровский товар"
    36 = "купецкий товар"
    37 = "галантерейный товар"
    38 = "заморский товар"
    39 = "меховой товар"
    40 = "рукодельный товар"
    41 = "домовый товар"
    42 = "надлежащий товар"
    43 = "харчевой припасы"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 1).Value = $updates[$row]
}
